$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before GB (old "nom" column), shifting nom/url_produit right to GC/GD.
$ws.Range("GB1").EntireColumn.Insert()

# Set header for the newly inserted column (new scrape timestamp).
$ws.Range("GB1").Value = "2026-02-05 12:32:36"

# Carry forward the last known price (previously in column GA) into the new column for rows 2-80.
for ($r = 2; $r -le 80; $r++) {
    $val = $ws.Range("GA$r").Value()
    $ws.Range("GB$r").Value = $val
}
